# Update "想去人数" (number of people wanting to go) counts by +1
# for the two events whose F2/F5 cells hold 340 and 288 respectively.
# This change needs to be applied on both the "展览" sheet and the
# "全部类型" sheet, which mirror the same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 341
    $ws.Range("F5").Value = 289
}
